$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest coin prices / 1h volume deltas pulled this run.
# Row = sheet row, Price = new column D text, Volume = new column E text
# ($null entries mean only the Volume% moved this refresh).
$updates = @(
    @{ Row = 2; Price = '51.438.49'; Volume = '  +0.32%  ' },
    @{ Row = 3; Price = '2.976.39'; Volume = '  +1.85%  ' },
    @{ Row = 4; Price = $null; Volume = '  -0.12%  ' },
    @{ Row = 5; Price = '378.43'; Volume = '  +2.29%  ' },
    @{ Row = 6; Price = '104.70'; Volume = '  +0.37%  ' },
    @{ Row = 7; Price = '0.540'; Volume = '  -0.08%  ' },
    @{ Row = 8; Price = $null; Volume = '  -0.09%  ' },
    @{ Row = 9; Price = '0.593'; Volume = '  +0.89%  ' },
    @{ Row = 10; Price = '37.19'; Volume = '  +1.08%  ' },
    @{ Row = 11; Price = '0.140'; Volume = '  +0.00%  ' },
    @{ Row = 12; Price = '0.0844'; Volume = '  +1.00%  ' },
    @{ Row = 13; Price = '3.438.49'; Volume = '  +1.36%  ' },
    @{ Row = 14; Price = '18.38'; Volume = '  +0.17%  ' },
    @{ Row = 15; Price = '7.59'; Volume = '  +2.62%  ' },
    @{ Row = 16; Price = '2.968.25'; Volume = '  +1.55%  ' },
    @{ Row = 17; Price = '0.965'; Volume = '  +2.53%  ' },
    @{ Row = 18; Price = '51.415.83'; Volume = '  +0.30%  ' },
    @{ Row = 19; Price = '3.33'; Volume = '  +2.57%  ' },
    @{ Row = 20; Price = '7.42'; Volume = '  +2.66%  ' },
    @{ Row = 21; Price = '12.92'; Volume = '  -0.27%  ' },
    @{ Row = 22; Price = '0.0₃0962'; Volume = '  +1.94%  ' },
    @{ Row = 23; Price = '69.31'; Volume = '  +1.36%  ' },
    @{ Row = 24; Price = '261.75'; Volume = '  +0.62%  ' },
    @{ Row = 25; Price = '2.83'; Volume = '  +4.97%  ' },
    @{ Row = 26; Price = '8.21'; Volume = '  +16.48%  ' },
    @{ Row = 27; Price = '7.63'; Volume = '  +21.59%  ' },
    @{ Row = 28; Price = $null; Volume = '  -2.12%  ' },
    @{ Row = 29; Price = '0.115'; Volume = '  +11.09%  ' },
    @{ Row = 30; Price = $null; Volume = '  +0.01%  ' },
    @{ Row = 31; Price = '25.88'; Volume = '  +0.33%  ' },
    @{ Row = 32; Price = '9.87'; Volume = '  -0.49%  ' },
    @{ Row = 33; Price = '35.09'; Volume = '  +0.70%  ' },
    @{ Row = 34; Price = $null; Volume = '  -2.10%  ' },
    @{ Row = 35; Price = '51.12'; Volume = '  +0.35%  ' },
    @{ Row = 36; Price = '0.0445'; Volume = '  +5.02%  ' },
    @{ Row = 37; Price = $null; Volume = '  +0.04%  ' },
    @{ Row = 38; Price = '3.02'; Volume = '  -0.27%  ' },
    @{ Row = 39; Price = '17.17'; Volume = '  +0.19%  ' },
    @{ Row = 40; Price = $null; Volume = '  -3.20%  ' },
    @{ Row = 41; Price = '1.85'; Volume = '  +0.06%  ' },
    @{ Row = 42; Price = '0.116'; Volume = '  +2.34%  ' },
    @{ Row = 43; Price = '125.36'; Volume = '  +5.12%  ' },
    @{ Row = 44; Price = '21.85'; Volume = '  -2.03%  ' },
    @{ Row = 45; Price = '0.284'; Volume = '  +16.42%  ' },
    @{ Row = 46; Price = $null; Volume = '  -1.28%  ' },
    @{ Row = 47; Price = $null; Volume = '  +3.06%  ' },
    @{ Row = 48; Price = '2.032.94'; Volume = '  +0.43%  ' },
    @{ Row = 49; Price = '3.23'; Volume = '  +1.49%  ' },
    @{ Row = 50; Price = $null; Volume = '  +8.15%  ' },
    @{ Row = 51; Price = '58.29'; Volume = '  +2.26%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $cell = $ws.Cells.Item($u.Row, 4)   # column D = Price

        if ($u.Price -match '^-?\d+(\.\d+)?$') {
            # Force text storage so e.g. "378.43" doesn't get reinterpreted
            # as a number the way Excel normally would on a General cell.
            $cell.NumberFormat = "@"
            $cell.Value = $u.Price
            $cell.Style = "Normal"
        } else {
            # Already unambiguous as text (multi-dot thousands grouping, etc.)
            $cell.Value = $u.Price
        }
    }

    $ws.Cells.Item($u.Row, 5).Value = $u.Volume   # column E = Volume(1h)
}
